$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Slide 6 (graphic-frame table "Sources of Finance") switches its
#    table style from the custom Table_0 style to a different built-in
#    style ({00428D80-AA62-4F06-A3CE-F6B338714C39}).
# ------------------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{00428D80-AA62-4F06-A3CE-F6B338714C39}")
    }
}

# ------------------------------------------------------------------
# 2) The deck's two embedded themes are swapped: the theme that is
#    actually applied to the slide master / slides (currently the
#    "Integral" palette) becomes the stock Office palette that used to
#    live in the unused (notes-master-only) theme part, and vice
#    versa. The live theme is reachable through any slide's
#    ThemeColorScheme, so recolor it in place to the Office values.
# ------------------------------------------------------------------
$officeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [math]::Floor($hex / 0x10000) % 0x100
    $g = [math]::Floor($hex / 0x100) % 0x100
    $b = $hex % 0x100
    $comRgb = ($b * 0x10000) + ($g * 0x100) + $r
    $themeColors.Colors($i).RGB = $comRgb
}
